$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data rows got re-sorted: the match that used to be stored on
# row 17 is now stored on row 18 (and vice versa), and likewise for rows
# 81/82. The "id" column (A) keeps its sequential value, but everything
# else describing the match (columns B:AC) moves with the match record.
# Swap the B:AC payload between each pair of rows.

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA" + ":AC$rowA")
    $rangeB = $ws.Range("B$rowB" + ":AC$rowB")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

Swap-Rows 17 18
Swap-Rows 81 82
